# Add pizza-type columns (Cheese / Pepperoni / Hawaiian / Meat Lovers) to the
# transactions sheet, splitting the old single "Pizza Type" column into four
# quantity columns and moving "Total" to the end.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ------------------------------------------------
$ws.Range("D1").Value = "Cheese Pizzas"
$ws.Range("E1").Value = "Pepperoni Pizzas"
$ws.Range("F1").Value = "Hawaiian Pizzas"
$ws.Range("G1").Value = "Meat Lovers Pizzas"
$ws.Range("H1").Value = "Total"

# Match the new header cells' formatting to the existing header style.
$ws.Range("D1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data row (row 2) ----------------------------------------------------
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 3
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 149.0775

# Match the new/ reused data cells' formatting to the existing data style.
$ws.Range("E2").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$ws.Range("F2:H2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Column widths -----------------------------------------------------
# (target stored widths, minus the 5/6-char Excel cell-padding constant that
# the ColumnWidth COM property re-adds internally)
$ws.Columns.Item(4).ColumnWidth = 21.63 - 5/6
$ws.Columns.Item(5).ColumnWidth = 20.13 - 5/6
$ws.Columns.Item(6).ColumnWidth = 23.25 - 5/6
$ws.Columns.Item(7).ColumnWidth = 31.88 - 5/6
$ws.Columns.Item(8).ColumnWidth = 20.63 - 5/6
